# Auto-generated edit script
# Applies cell-value corrections to the Excalibur_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1821.5555
$ws.Range("J17").Value = 1899.25
$ws.Range("L17").Value = 5697.75
$ws.Range("N17").Value = -6033.75
# Row 19
$ws.Range("H19").Value = 325.42856
$ws.Range("I19").Value = 292
$ws.Range("J19").Value = 350.5
$ws.Range("K19").Value = 292
$ws.Range("L19").Value = 350.5
$ws.Range("M19").Value = -117
$ws.Range("N19").Value = -700.5
# Row 33
$ws.Range("H33").Value = 825.44446
$ws.Range("I33").Value = 754.8333
$ws.Range("J33").Value = 966.6667
$ws.Range("K33").Value = 754.8333
$ws.Range("L33").Value = 966.6667
$ws.Range("M33").Value = -525.8333
$ws.Range("N33").Value = -1424.6667
# Row 41
$ws.Range("H41").Value = 83895.164
$ws.Range("I41").Value = 142.25
$ws.Range("K41").Value = 142.25
$ws.Range("M41").Value = 297.75
# Row 62
$ws.Range("H62").Value = 25535.533
$ws.Range("I62").Value = 30379.125
$ws.Range("K62").Value = 30379.125
$ws.Range("M62").Value = -29755.125
# Row 65
$ws.Range("H65").Value = 25535.533
$ws.Range("I65").Value = 30379.125
$ws.Range("K65").Value = 151895.625
$ws.Range("M65").Value = -148775.625
# Row 113
$ws.Range("H113").Value = 2998.75
$ws.Range("I113").Value = 2998.75
$ws.Range("K113").Value = 2998.75
$ws.Range("M113").Value = 255.25
# Row 135
$ws.Range("H135").Value = 2168.2942
$ws.Range("I135").Value = 2239.5
$ws.Range("K135").Value = 20155.5
$ws.Range("M135").Value = -17620.5
# Row 138
$ws.Range("H138").Value = 3029.9355
$ws.Range("J138").Value = 5611.1113
$ws.Range("L138").Value = 16833.3339
$ws.Range("N138").Value = -27113.3339

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3724514.8
$ws.Range("I61").Value = 3724514.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3724514.8
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = -3724302.8
# Row 81
$ws.Range("H81").Value = 80181
$ws.Range("J81").Value = 80181
$ws.Range("L81").Value = 80181
$ws.Range("N81").Value = -82177
# Row 84
$ws.Range("H84").Value = 80181
$ws.Range("J84").Value = 80181
$ws.Range("L84").Value = 240543
$ws.Range("N84").Value = -250527
# Row 132
$ws.Range("H132").Value = 508025.28
$ws.Range("I132").Value = 555318.8
$ws.Range("K132").Value = 1665956.4
$ws.Range("M132").Value = -1663426.4
# Row 136
$ws.Range("H136").Value = 3724514.8
$ws.Range("I136").Value = 3724514.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11173544.4
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -11170994.4

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1513.5333
$ws.Range("I20").Value = 1626.15
$ws.Range("J20").Value = 1288.3
$ws.Range("K20").Value = 1626.15
$ws.Range("L20").Value = 1288.3
$ws.Range("M20").Value = -1379.15
$ws.Range("N20").Value = -1782.3
# Row 132
$ws.Range("H132").Value = 99970
$ws.Range("J132").Value = 99970
$ws.Range("L132").Value = 99970
$ws.Range("N132").Value = -110090

$ws = $wb.Worksheets.Item("CRP")
# Row 47
$ws.Range("H47").Value = 59500
$ws.Range("J47").Value = 59500
$ws.Range("L47").Value = 59500
$ws.Range("N47").Value = -60632
# Row 105
$ws.Range("H105").Value = 69948.39999999999
$ws.Range("I105").Value = 87060.5
$ws.Range("K105").Value = 87060.5
$ws.Range("M105").Value = -85313.5
# Row 132
$ws.Range("H132").Value = 8634726
$ws.Range("I132").Value = 15600.808
$ws.Range("J132").Value = 83333810
$ws.Range("K132").Value = 46802.424
$ws.Range("L132").Value = 250001430
$ws.Range("M132").Value = -44272.424
$ws.Range("N132").Value = -250006490
# Row 134
$ws.Range("H134").Value = 1753.4445
$ws.Range("I134").Value = 1753.4445
$ws.Range("K134").Value = 5260.333500000001
$ws.Range("M134").Value = -2725.333500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 800
$ws.Range("J98").Value = 800
$ws.Range("L98").Value = 2400
$ws.Range("N98").Value = -5396
# Row 114
$ws.Range("H114").Value = 1466.36
$ws.Range("I114").Value = 126.9375
$ws.Range("K114").Value = 380.8125
$ws.Range("M114").Value = 2873.1875

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6009.3
$ws.Range("I70").Value = 6870.5713
$ws.Range("K70").Value = 6870.5713
$ws.Range("M70").Value = -6600.5713
# Row 73
$ws.Range("H73").Value = 6009.3
$ws.Range("I73").Value = 6870.5713
$ws.Range("K73").Value = 6870.5713
$ws.Range("M73").Value = -5934.5713
# Row 97
$ws.Range("H97").Value = 10952.167
$ws.Range("J97").Value = 19927
$ws.Range("L97").Value = 19927
$ws.Range("N97").Value = -20919
# Row 102
$ws.Range("H102").Value = 3016.423
$ws.Range("I102").Value = 1706
$ws.Range("K102").Value = 1706
$ws.Range("M102").Value = -84
# Row 114
$ws.Range("H114").Value = 80722
$ws.Range("J114").Value = 80722
$ws.Range("L114").Value = 80722
$ws.Range("N114").Value = -89400
# Row 132
$ws.Range("H132").Value = 862465.4
$ws.Range("I132").Value = 1005385.2
$ws.Range("J132").Value = 4946.5
$ws.Range("K132").Value = 3016155.6
$ws.Range("L132").Value = 14839.5
$ws.Range("M132").Value = -3013625.6
$ws.Range("N132").Value = -19899.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 49694.24
$ws.Range("I22").Value = 78145.234
$ws.Range("J22").Value = 3461.375
$ws.Range("K22").Value = 78145.234
$ws.Range("L22").Value = 3461.375
$ws.Range("M22").Value = -77850.234
$ws.Range("N22").Value = -4051.375
# Row 27
$ws.Range("H27").Value = 49694.24
$ws.Range("I27").Value = 78145.234
$ws.Range("J27").Value = 3461.375
$ws.Range("K27").Value = 78145.234
$ws.Range("L27").Value = 3461.375
$ws.Range("M27").Value = -78038.234
$ws.Range("N27").Value = -3675.375
# Row 45
$ws.Range("H45").Value = 15000
$ws.Range("I45").Value = 15000
$ws.Range("K45").Value = 15000
$ws.Range("M45").Value = -14593
# Row 46
$ws.Range("H46").Value = 1184.1538
$ws.Range("I46").Value = 1163.091
$ws.Range("J46").Value = 1199.6
$ws.Range("K46").Value = 1163.091
$ws.Range("L46").Value = 1199.6
$ws.Range("M46").Value = -975.0909999999999
$ws.Range("N46").Value = -1575.6
# Row 55
$ws.Range("H55").Value = 1458.625
$ws.Range("J55").Value = 1886.3334
$ws.Range("L55").Value = 1886.3334
$ws.Range("N55").Value = -2232.3334
# Row 61
$ws.Range("H61").Value = 2532.5715
$ws.Range("I61").Value = 1208.3636
$ws.Range("J61").Value = 3989.2
$ws.Range("K61").Value = 1208.3636
$ws.Range("L61").Value = 3989.2
$ws.Range("M61").Value = -1006.3636
$ws.Range("N61").Value = -4393.2
# Row 103
$ws.Range("H103").Value = 116665.336
$ws.Range("J103").Value = 116665.336
$ws.Range("L103").Value = 116665.336
$ws.Range("N103").Value = -119009.336
# Row 113
$ws.Range("H113").Value = 2532.5715
$ws.Range("I113").Value = 1208.3636
$ws.Range("J113").Value = 3989.2
$ws.Range("K113").Value = 1208.3636
$ws.Range("L113").Value = 3989.2
$ws.Range("M113").Value = 961.6364000000001
$ws.Range("N113").Value = -8329.200000000001
# Row 132
$ws.Range("H132").Value = 667782.4399999999
$ws.Range("I132").Value = 770737.5600000001
$ws.Range("K132").Value = 2312212.68
$ws.Range("M132").Value = -2309682.68

$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Range("H123").Value = 52737.5
$ws.Range("J123").Value = 52737.5
$ws.Range("L123").Value = 52737.5
$ws.Range("N123").Value = -62537.5
# Row 132
$ws.Range("H132").Value = 4378077
$ws.Range("I132").Value = 5592599
$ws.Range("J132").Value = 5798.6
$ws.Range("K132").Value = 16777797
$ws.Range("L132").Value = 17395.8
$ws.Range("M132").Value = -16775267
$ws.Range("N132").Value = -22455.8
# Row 136
$ws.Range("H136").Value = 10303517
$ws.Range("I136").Value = 12296651
$ws.Range("J136").Value = 5658.1665
$ws.Range("K136").Value = 36889953
$ws.Range("L136").Value = 16974.4995
$ws.Range("M136").Value = -36887403
$ws.Range("N136").Value = -22074.4995
